$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds for row 2 (The Strongest vs GV San Jose) ---
$ws.Range("G2").Value = 1.36
$ws.Range("K2").Value = 2.75
$ws.Range("L2").Value = 6.5
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 1.48
$ws.Range("R2").Value = 2.6
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 8
$ws.Range("Y2").Value = 9
$ws.Range("AC2").Value = 19
$ws.Range("AE2").Value = 17
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AJ2").Value = 21
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 3.6
$ws.Range("AS2").Value = 81

# --- Update odds for row 4 (Guarani vs Novorizontino) ---
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 2.8
$ws.Range("L4").Value = 3.75
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("T4").Value = 2.08
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 12
$ws.Range("AH4").Value = 6
$ws.Range("AL4").Value = 34
$ws.Range("AO4").Value = 19
$ws.Range("AW4").Value = 4.5
$ws.Range("AX4").Value = 19
$ws.Range("AY4").Value = 41

# --- Remove the Dep. Pasto vs Aguilas match (old row 5); remaining rows shift up ---
$ws.Rows("5").Delete()
